$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 482.66666
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 482.66666
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 482.66666
$ws.Range("N12").Value = -822.66666
$ws.Range("M12").ClearContents()

$ws.Range("H42").Value = 441.5
$ws.Range("I42").Value = 441.5
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 1324.5
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -1094.5
$ws.Range("N42").ClearContents()

$ws.Range("H43").Value = 8000
$ws.Range("J43").Value = 8000
$ws.Range("L43").Value = 8000
$ws.Range("N43").Value = -8138

$ws.Range("H51").Value = 14599.4
$ws.Range("J51").Value = 14999.5
$ws.Range("L51").Value = 14999.5
$ws.Range("N51").Value = -15967.5

$ws.Range("H138").Value = 1518.3158
$ws.Range("J138").Value = 1512.8572
$ws.Range("L138").Value = 4538.571599999999
$ws.Range("N138").Value = -14818.5716

$ws.Range("H141").Value = 13138.4
$ws.Range("I141").Value = 15364
$ws.Range("K141").Value = 46092
$ws.Range("M141").Value = -40912
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1786.0834
$ws.Range("J88").Value = 2061.8572
$ws.Range("L88").Value = 2061.8572
$ws.Range("N88").Value = -2873.8572

$ws.Range("H91").Value = 1786.0834
$ws.Range("J91").Value = 2061.8572
$ws.Range("L91").Value = 2061.8572
$ws.Range("N91").Value = -4869.8572

$ws.Range("H110").Value = 6314.769
$ws.Range("I110").Value = 7048.143
$ws.Range("K110").Value = 7048.143
$ws.Range("M110").Value = -5003.143
$ws.Range("N110").ClearContents()

$ws.Range("H122").Value = 2378.6
$ws.Range("I122").Value = 2297.6667
$ws.Range("K122").Value = 6893.000100000001
$ws.Range("M122").Value = -4443.000100000001
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 3133.5
$ws.Range("I132").Value = 3133.5
$ws.Range("K132").Value = 9400.5
$ws.Range("M132").Value = -6870.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 8605.762000000001
$ws.Range("I86").Value = 10387.538
$ws.Range("J86").Value = 5710.375
$ws.Range("K86").Value = 10387.538
$ws.Range("L86").Value = 5710.375
$ws.Range("M86").Value = -9264.538
$ws.Range("N86").Value = -7956.375

$ws.Range("H89").Value = 8605.762000000001
$ws.Range("I89").Value = 10387.538
$ws.Range("J89").Value = 5710.375
$ws.Range("K89").Value = 51937.69
$ws.Range("L89").Value = 28551.875
$ws.Range("M89").Value = -46321.69
$ws.Range("N89").Value = -39783.875

$ws.Range("H134").Value = 1552.8182
$ws.Range("I134").Value = 1552.8182
$ws.Range("K134").Value = 4658.4546
$ws.Range("M134").Value = -2123.4546

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 420
$ws.Range("I7").Value = 155.6
$ws.Range("J7").Value = 750.5
$ws.Range("K7").Value = 155.6
$ws.Range("L7").Value = 750.5
$ws.Range("M7").Value = -42.59999999999999
$ws.Range("N7").Value = -976.5

$ws.Range("H31").Value = 3460.625
$ws.Range("I31").Value = 2496.25
$ws.Range("K31").Value = 2496.25
$ws.Range("M31").Value = -2201.25
$ws.Range("N31").ClearContents()

$ws.Range("H34").Value = 3460.625
$ws.Range("I34").Value = 2496.25
$ws.Range("K34").Value = 2496.25
$ws.Range("M34").Value = -2294.25
$ws.Range("N34").ClearContents()

$ws.Range("H69").Value = 14544.75
$ws.Range("I69").Value = 14544.75
$ws.Range("K69").Value = 14544.75
$ws.Range("M69").Value = -13795.75

$ws.Range("H72").Value = 14544.75
$ws.Range("I72").Value = 14544.75
$ws.Range("K72").Value = 43634.25
$ws.Range("M72").Value = -39890.25

$ws.Range("H107").Value = 789.5294
$ws.Range("I107").Value = 659.9231
$ws.Range("J107").Value = 1210.75
$ws.Range("K107").Value = 659.9231
$ws.Range("L107").Value = 1210.75
$ws.Range("M107").Value = 1260.0769
$ws.Range("N107").Value = -5050.75

$ws.Range("H122").Value = 1962
$ws.Range("I122").Value = 1962
$ws.Range("K122").Value = 5886
$ws.Range("M122").Value = -3436

$ws.Range("H134").Value = 2110.5
$ws.Range("I134").Value = 1862.5
$ws.Range("K134").Value = 5587.5
$ws.Range("M134").Value = -3052.5
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 4066.6667
$ws.Range("J86").Value = 7000
$ws.Range("L86").Value = 21000
$ws.Range("N86").Value = -23372

$ws.Range("H89").Value = 4066.6667
$ws.Range("J89").Value = 7000
$ws.Range("L89").Value = 63000
$ws.Range("N89").Value = -74856

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3071.75
$ws.Range("I80").Value = 3295.6667
$ws.Range("K80").Value = 3295.6667
$ws.Range("M80").Value = -2297.6667
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 3071.75
$ws.Range("I83").Value = 3295.6667
$ws.Range("K83").Value = 16478.3335
$ws.Range("M83").Value = -11486.3335
$ws.Range("N83").ClearContents()

$ws.Range("H102").Value = 1135.9
$ws.Range("I102").Value = 884.8333
$ws.Range("J102").Value = 1512.5
$ws.Range("K102").Value = 884.8333
$ws.Range("L102").Value = 1512.5
$ws.Range("M102").Value = 737.1667
$ws.Range("N102").Value = -4756.5

$ws.Range("H122").Value = 2912.25
$ws.Range("I122").Value = 2549.6667
$ws.Range("K122").Value = 7649.000100000001
$ws.Range("M122").Value = -5199.000100000001
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 3023.8125
$ws.Range("I132").Value = 2959.7693
$ws.Range("J132").Value = 3301.3333
$ws.Range("K132").Value = 8879.3079
$ws.Range("L132").Value = 9903.999899999999
$ws.Range("M132").Value = -6349.3079
$ws.Range("N132").Value = -14963.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1946
$ws.Range("I7").Value = 1946
$ws.Range("K7").Value = 1946
$ws.Range("M7").Value = -1834

$ws.Range("H40").Value = 3235.6
$ws.Range("I40").Value = 2479.5715
$ws.Range("K40").Value = 2479.5715
$ws.Range("M40").Value = -2343.5715
$ws.Range("N40").ClearContents()

$ws.Range("H46").Value = 1524.125
$ws.Range("I46").Value = 1107
$ws.Range("J46").Value = 3331.6667
$ws.Range("K46").Value = 1107
$ws.Range("L46").Value = 3331.6667
$ws.Range("M46").Value = -919
$ws.Range("N46").Value = -3707.6667

$ws.Range("H68").Value = 4250
$ws.Range("I68").Value = 3500
$ws.Range("J68").Value = 5000
$ws.Range("K68").Value = 3500
$ws.Range("L68").Value = 5000
$ws.Range("M68").Value = -2751
$ws.Range("N68").Value = -6498

$ws.Range("H71").Value = 4250
$ws.Range("I71").Value = 3500
$ws.Range("J71").Value = 5000
$ws.Range("K71").Value = 17500
$ws.Range("L71").Value = 25000
$ws.Range("M71").Value = -13756
$ws.Range("N71").Value = -32488

$ws.Range("H122").Value = 3619
$ws.Range("I122").Value = 3526
$ws.Range("K122").Value = 10578
$ws.Range("M122").Value = -8128
$ws.Range("N122").ClearContents()

$ws.Range("H126").Value = 1946
$ws.Range("I126").Value = 1946
$ws.Range("K126").Value = 5838
$ws.Range("M126").Value = -3368

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9998.5
$ws.Range("J62").Value = 9998.5
$ws.Range("L62").Value = 9998.5
$ws.Range("N62").Value = -11246.5

$ws.Range("H65").Value = 9998.5
$ws.Range("J65").Value = 9998.5
$ws.Range("L65").Value = 49992.5
$ws.Range("N65").Value = -56232.5

$ws.Range("H122").Value = 3342.7273
$ws.Range("I122").Value = 3177
$ws.Range("K122").Value = 9531
$ws.Range("M122").Value = -7081
$ws.Range("N122").ClearContents()

$ws.Range("H126").Value = 4420.9165
$ws.Range("I126").Value = 4167.9414
$ws.Range("K126").Value = 12503.8242
$ws.Range("M126").Value = -10033.8242
$ws.Range("N126").ClearContents()

$ws.Range("H132").Value = 3009.3684
$ws.Range("I132").Value = 2481.0588
$ws.Range("K132").Value = 7443.176399999999
$ws.Range("M132").Value = -4913.176399999999
$ws.Range("N132").ClearContents()
